# Regenerate merged AHB files
# - Rename the "_old" / "_new" suffixed header labels in row 1 to
#   "_FV2304" / "_FV2310" respectively.
# - Wrap the used range A1:U57 in an Excel Table ("Table1").
# - Freeze the header row (row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:U1) -----------------------------------------
$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the data range into a proper Excel Table ----------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
